# Apply BOM correction: fix mouser part number text in cell F6
# (80-C0805C226M8 -> 80-C0805C226M8P, trailing space included)
# and restore the active selection to C7 (as saved by the author).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = "80-C0805C226M8P "

$ws.Range("C7").Select()
